# Update cryptocurrency price/volume snapshot data on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.930.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5162"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3716"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07188"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8984"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.857.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.254"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008490"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.956.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.027"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.102.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.427"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.781"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.899"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.745"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09175"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05028"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7533"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.991"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.171"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.269"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01992"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5566"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.486"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.071"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.575"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.724"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "115.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1498"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4768"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.562"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
